{"js": "// Update the date and all the three-digit-by-one-digit multiplication\n// problems in the document. Every `<w:t>` run in this worksheet is unique,\n// so a literal, case-sensitive search-and-replace for each old value is\n// unambiguous.\nconst replacements = [\n  [\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"],\n  [\"211\u00d79=\", \"624\u00d78=\"],\n  [\"800\u00d77=\", \"161\u00d79=\"],\n  [\"622\u00d75=\", \"628\u00d77=\"],\n  [\"665\u00d73=\", \"341\u00d78=\"],\n  [\"796\u00d75=\", \"980\u00d79=\"],\n  [\"121\u00d73=\", \"945\u00d75=\"],\n  [\"803\u00d77=\", \"275\u00d78=\"],\n  [\"503\u00d75=\", \"881\u00d75=\"],\n  [\"951\u00d76=\", \"612\u00d74=\"],\n  [\"558\u00d75=\", \"403\u00d79=\"],\n  [\"230\u00d77=\", \"914\u00d79=\"],\n  [\"464\u00d73=\", \"430\u00d77=\"],\n  [\"955\u00d74=\", \"293\u00d74=\"],\n  [\"804\u00d79=\", \"743\u00d75=\"],\n  [\"427\u00d76=\", \"431\u00d76=\"],\n  [\"431\u00d79=\", \"333\u00d79=\"],\n  [\"208\u00d77=\", \"918\u00d79=\"],\n  [\"609\u00d74=\", \"925\u00d78=\"],\n  [\"995\u00d78=\", \"690\u00d72=\"],\n  [\"785\u00d78=\", \"791\u00d78=\"],\n  [\"781\u00d73=\", \"784\u00d73=\"],\n  [\"746\u00d79=\", \"548\u00d75=\"],\n  [\"293\u00d72=\", \"699\u00d74=\"],\n  [\"767\u00d77=\", \"193\u00d73=\"],\n  [\"577\u00d73=\", \"630\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and all the three-digit-by-one-digit multiplication\n# problems in the document. Every run of text in this worksheet is unique,\n# so a literal Find/Replace for each old value is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"),\n    @(\"211\u00d79=\", \"624\u00d78=\"),\n    @(\"800\u00d77=\", \"161\u00d79=\"),\n    @(\"622\u00d75=\", \"628\u00d77=\"),\n    @(\"665\u00d73=\", \"341\u00d78=\"),\n    @(\"796\u00d75=\", \"980\u00d79=\"),\n    @(\"121\u00d73=\", \"945\u00d75=\"),\n    @(\"803\u00d77=\", \"275\u00d78=\"),\n    @(\"503\u00d75=\", \"881\u00d75=\"),\n    @(\"951\u00d76=\", \"612\u00d74=\"),\n    @(\"558\u00d75=\", \"403\u00d79=\"),\n    @(\"230\u00d77=\", \"914\u00d79=\"),\n    @(\"464\u00d73=\", \"430\u00d77=\"),\n    @(\"955\u00d74=\", \"293\u00d74=\"),\n    @(\"804\u00d79=\", \"743\u00d75=\"),\n    @(\"427\u00d76=\", \"431\u00d76=\"),\n    @(\"431\u00d79=\", \"333\u00d79=\"),\n    @(\"208\u00d77=\", \"918\u00d79=\"),\n    @(\"609\u00d74=\", \"925\u00d78=\"),\n    @(\"995\u00d78=\", \"690\u00d72=\"),\n    @(\"785\u00d78=\", \"791\u00d78=\"),\n    @(\"781\u00d73=\", \"784\u00d73=\"),\n    @(\"746\u00d79=\", \"548\u00d75=\"),\n    @(\"293\u00d72=\", \"699\u00d74=\"),\n    @(\"767\u00d77=\", \"193\u00d73=\"),\n    @(\"577\u00d73=\", \"630\u00d72=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
